# "Diseño de pruebas + AVL TAD terminado"
#
# The first table's first cell currently holds a single bold run
# "TAD ". We need to complete the title by adding a second bold run
# "AVL Tree" right after it (and before the existing _GoBack bookmark),
# producing the visible text "TAD AVL Tree" as two runs:
#   <w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">TAD </w:t></w:r>
#   <w:r><w:rPr><w:b/></w:rPr><w:t>AVL Tree</w:t></w:r>

$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$cell  = $table.Cell(1, 1)
$para  = $cell.Range.Paragraphs.Item(1)

# Collapse to the end of the paragraph's text, excluding the paragraph
# mark itself (so the bookmark that follows stays right after our text).
$insertionPoint = $para.Range.Duplicate
$insertionPoint.Collapse(0)              # wdCollapseEnd
$insertionPoint.MoveEnd(1, -1) | Out-Null
$startPos = $insertionPoint.Start

# Append the new text; at this point it is merged into the existing
# "TAD " run (same bold formatting), carrying that run's text onward.
$insertionPoint.InsertAfter("AVL Tree")

# Force the appended text to live in its own run (matching the target
# markup, which has "TAD " and "AVL Tree" as two separate <w:r>
# elements) by round-tripping Bold off/on purely over that new span.
$newText = $d.Range($startPos, $startPos + 8)
$newText.Bold = 0
$newText = $d.Range($startPos, $startPos + 8)
$newText.Bold = 1
